$d = $word.ActiveDocument

# Anchor text: the paragraph that the three new bullet points must be inserted before.
$anchorText = "Developed and deployed custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering"

$newBullets = @(
    "• Developed meta-analytical techniques that identified systematic data quality issues across 20+ years of voter registration data",
    "• Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters",
    "• Created fraud detection systems analyzing 5+ terabyte datasets, uncovering demographic miscoding patterns across 2,000+ precincts"
)

foreach ($bulletText in $newBullets) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$anchorText*") {
            $r = $p.Range
            $startPos = $r.Start
            $r.InsertParagraphBefore()
            $newPara = $d.Range($startPos, $startPos)
            $newPara.InsertAfter($bulletText)
            break
        }
    }
}

# Remove the old "Created fraud detection systems for campaign finance..." bullet
# (it followed the "...170% more viable targets" paragraph, under the same PARTNER role).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Created fraud detection systems for campaign finance data analysis across multi-terabyte datasets*") {
        $p.Range.Delete()
        break
    }
}
